# Correct status name(s) in the "liste essais cliniques identifies" sheet.
#
# Changes needed (matching the shared-string diff):
#   "bleu"                                                  -> "noir"
#   "résultat et / ou publication posté"                    -> "résultat postés ou publiés"
#   "pas de résultat ni de publication"                     -> "pas de résultat postés ni publiés"
#   "résultat et / ou publication posté dans les 12 mois"   -> "résultat postés ou publiés dans les 12 mois"
#   "résultat et / ou publication posté dans les 36 mois"   -> "résultat postés ou publiés dans les 36 mois"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "bleu"                                                = "noir"
    "résultat et / ou publication posté"                  = "résultat postés ou publiés"
    "pas de résultat ni de publication"                   = "pas de résultat postés ni publiés"
    "résultat et / ou publication posté dans les 12 mois" = "résultat postés ou publiés dans les 12 mois"
    "résultat et / ou publication posté dans les 36 mois" = "résultat postés ou publiés dans les 36 mois"
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $usedRange.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($null -ne $val -and $replacements.ContainsKey($val)) {
            $cell.Value = $replacements[$val]
        }
    }
}
